$d = $word.ActiveDocument

# --- 1. Update the first paragraph: keep the existing sentence (now with two
#        trailing spaces) and append a new, dark-red run with the
#        "branch alternate" remark. ---
$p1 = $d.Paragraphs(1)
$insertPos = $p1.Range.Start + $p1.Range.Text.TrimEnd([char]13, [char]7).Length

$spaceRange = $d.Range($insertPos, $insertPos)
$spaceRange.InsertAfter("  ")

$notePos = $spaceRange.End
$noteRange = $d.Range($notePos, $notePos)
$noteText = "(This is a change " + [char]0x2013 + " Version for branch alternate)"
$noteRange.InsertAfter($noteText)
$noteRange.Font.Color = 192

# --- 2. Mark the "Normal (Web)" style as semi-hidden (hidden until used). ---
$style = $d.Styles("Normal (Web)")
$style.Visibility = $false
